$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.028.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.06%  '

$ws.Range("D3").Value = "'3.393.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.11%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'559.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.62%  '

$ws.Range("D6").Value = "'174.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.89%  '

$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.78%  '

$ws.Range("D8").Value = "'3.383.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.06%  '

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +9.99%  '

$ws.Range("D11").Value = "'0.633"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.23%  '

$ws.Range("D12").Value = "'54.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.03%  '

$ws.Range("D13").Value = "'0.0000277"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.18%  '

$ws.Range("D14").Value = "'9.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.56%  '

$ws.Range("D15").Value = "'3.917.95"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.98%  '

$ws.Range("D16").Value = "'18.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("E17").Value = '  +3.20%  '

$ws.Range("D18").Value = "'3.372.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.42%  '

$ws.Range("D19").Value = "'64.908.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").Value = "'11.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.34%  '

$ws.Range("D21").Value = "'0.993"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.31%  '

$ws.Range("D22").Value = "'467.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +12.62%  '

$ws.Range("E23").Value = '  +12.66%  '

$ws.Range("E24").Value = '  +3.51%  '

$ws.Range("D25").Value = "'86.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.98%  '

$ws.Range("D26").Value = "'13.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.67%  '

$ws.Range("D27").Value = "'2.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.36%  '

$ws.Range("D28").Value = "'10.82"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.77%  '

$ws.Range("D29").Value = "'8.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.89%  '

$ws.Range("D30").Value = "'30.81"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.85%  '

$ws.Range("D31").Value = "'6.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +7.14%  '

$ws.Range("D32").Value = "'11.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.16%  '

$ws.Range("D33").Value = "'572.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.73%  '

$ws.Range("D34").Value = "'61.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.47%  '

$ws.Range("E35").Value = '  +2.68%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D37").Value = "'3.61"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.13%  '

$ws.Range("E38").Value = '  -3.56%  '

$ws.Range("D39").Value = "'35.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.16%  '

$ws.Range("D40").Value = "'0.0₃0742"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.23%  '

$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").Value = "'3.087.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").Value = "'2.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.86%  '

$ws.Range("D45").Value = "'0.0415"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +4.67%  '

$ws.Range("D46").Value = "'0.134"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.72%  '

$ws.Range("E47").Value = '  +2.38%  '

$ws.Range("D48").Value = "'3.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.28%  '

$ws.Range("D49").Value = "'2.59"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.75%  '

$ws.Range("D50").Value = "'138.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.68%  '

$ws.Range("D51").Value = "'8.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.90%  '
